# Schedule workbook edit: "changed event order and names"
#
# Summary of the change (Day 2 / 14-05-2025 afternoon block):
#  - The old "Advanced Quarto" slot (10:30-12:00) is split into two
#    half-hour slots: "Quarto (continuation)" (10:30-11:00) and a new
#    "Dynamic plotting" slot (11:00-12:00) with its own slide link.
#  - The old "Intro to dynamic plotting" slot (13:00-15:00) is split into
#    "Dynamic plotting (continuation)" (13:00-14:00) and "Quarto
#    Dashboards" (14:00-15:00, reusing the slide_quarto2 link).
#  - The old "Advanced dynamic plotting" slot becomes
#    "Quarto Dashboards (continuation)".
#  - A couple of Day-1 topic names are shortened, and the defined name
#    "schedule" grows from K23 to K25 to keep covering the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Shorten a couple of Day-1 topic names (same cells, new text).
# ---------------------------------------------------------------------
$ws.Range("E7").Value  = "Vector and bitmap"
$ws.Range("E10").Value = "Single-cell"

# ---------------------------------------------------------------------
# 2. Day-2 header: "Intro to Quarto" -> "Quarto"
# ---------------------------------------------------------------------
$ws.Range("E12").Value = "Quarto"

# ---------------------------------------------------------------------
# 3. Insert two new rows in the Day-2 afternoon block so the existing
#    rows shift down to make room for the new "Dynamic plotting" and
#    "Quarto Dashboards" split sessions. Inserting row 15 then row 18
#    (rather than both at once) lands every later row exactly where the
#    target layout needs it, since a fresh blank row is produced right
#    after each of the two sessions being split.
# ---------------------------------------------------------------------
$ws.Rows("15:15").Insert()
$ws.Rows("18:18").Insert()

# ---------------------------------------------------------------------
# 4. Row 14: "Advanced Quarto" -> "Quarto (continuation)", now only
#    10:30-11:00, taught by Katja Kozjek (LM, MR, LV), no slide link.
# ---------------------------------------------------------------------
$ws.Range("D14").Value = 0.458333333333333
$ws.Range("E14").Value = "Quarto (continuation)"
$ws.Range("F14").Value = "Katja Kozjek"
$ws.Range("G14").Value = "LM, MR, LV"
$ws.Range("H14").ClearContents()

# ---------------------------------------------------------------------
# 5. Row 15 (new): "Dynamic plotting", 11:00-12:00.
# ---------------------------------------------------------------------
$ws.Range("C15").Value = 0.458333333333333
$ws.Range("D15").Value = 0.5
$ws.Range("E15").Value = "Dynamic plotting"
$ws.Range("F15").Value = "Katja Kozjek"
$ws.Range("G15").Value = "LM, MR, LV"
$ws.Range("H15").Value = "topics/dynamic_plotting/slide_dynamic_plot1.html"

# Row 16 (Lunch, 12:00-13:00) already lines up correctly after the insert.

# ---------------------------------------------------------------------
# 6. Row 17: "Intro to dynamic plotting" -> "Dynamic plotting
#    (continuation)", now only 13:00-14:00, no slide link (it moved to
#    row 15 above).
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 0.583333333333333
$ws.Range("E17").Value = "Dynamic plotting (continuation)"
$ws.Range("F17").Value = "Katja Kozjek"
$ws.Range("G17").Value = "LM, MR, LV"
$ws.Range("H17").ClearContents()

# ---------------------------------------------------------------------
# 7. Row 18 (new): "Quarto Dashboards", 14:00-15:00, reusing the
#    slide_quarto2 link.
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 0.583333333333333
$ws.Range("D18").Value = 0.625
$ws.Range("E18").Value = "Quarto Dashboards"
$ws.Range("F18").Value = "Katja Kozjek"
$ws.Range("G18").Value = "LM, MR, LV"
$ws.Range("H18").Value = "topics/Quarto/slide_quarto2.html"

# Row 19 (Fika, 15:00-15:30... 13:30-15:30 slot) already lines up
# correctly (time values untouched by the insert).

# ---------------------------------------------------------------------
# 8. Row 20: "Advanced dynamic plotting" -> "Quarto Dashboards
#    (continuation)", taught by Katja Kozjek (LM, MR, LV) instead of
#    Lokesh Mano (KK, MR, LV).
# ---------------------------------------------------------------------
$ws.Range("E20").Value = "Quarto Dashboards (continuation)"
$ws.Range("F20").Value = "Katja Kozjek"
$ws.Range("G20").Value = "LM, MR, LV"

# Rows 21-31 (Wrap-up day II, Dinner, Day-3 R-Shiny schedule, ...) keep
# their original text/values; only their row numbers shifted, which the
# two row-inserts above already handled.

# ---------------------------------------------------------------------
# 9. Defined name "schedule" now covers the two extra rows.
# ---------------------------------------------------------------------
$nm = $wb.Names.Item("Sheet1!schedule")
$nm.RefersTo = "=Sheet1!`$A`$1:`$K`$25"

# ---------------------------------------------------------------------
# 10. Update the saved view/selection state.
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E20").Select()
